$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Main"

# Change every used cell's font to Arial 10 (was Aptos Narrow 11 / Aptos Narrow 11 variants).
# Cells are grouped by their pre-existing style so that each distinct original
# style only needs to be re-stamped once (keeps the produced style table tight,
# matching how Excel collapses identical formats to a single cellXf/font).
$group0  = @("A2","H2","H3","B4","H4","H5","H6","H7","H9","H10")
$group1  = @("A1")
$group2  = @("I2")
$group3  = @("I3","I4","I5","I6","I7")
$group4  = @("B5")
$group5  = @("J3","J5","J6")

foreach ($g in @($group0, $group1, $group2, $group3, $group4, $group5)) {
  foreach ($addr in $g) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Arial"
    $c.Font.Size = 10
  }
}

# Reset the view: scroll back to the top-left corner and move the selection.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("A11").Select()
